# "Add some deck related" — add a new ship (row 7), bump O3, update the
# style-comment on column G (style) to drop the obsolete "1/2 origin" note,
# and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- O3: size flag for the first ship becomes 1 -----------------------
$ws.Range("O3").Value = 1

# --- New row 7: "卡鲁提拉号" -------------------------------------------
# Copy the formatting from row 6 (B/C use the shared "s=1" 宋体 style) so
# the new row's name + size columns match the rest of the table.
$ws.Range("B6:C6").Copy()
$ws.Range("B7:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "卡鲁提拉号"
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 25
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 9
$ws.Range("L7").Value = 24
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 2
$ws.Range("P7").Value = 0

# --- Comment on E1: drop the obsolete "1/2 origin" explanation --------
$comment = $ws.Range("E1").Comment
$comment.Text("Yujie Liu:" + "`n" + "shipstyle" + "`n" + "`n")

# --- Move the selection to match the author's last edit position ------
$ws.Range("E5").Select()
